# Updated symbol list on Sat Jan 28 15:39:55 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) readings for
# the coins listed in the sheet. Values are stored as text (matching the
# original inline-string cells), so each target cell is explicitly formatted
# as Text before the new reading is written - this avoids Excel silently
# re-interpreting strings such as "0.9300" or "0.0001300" as numbers and
# dropping the significant trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "305.95"
    "E2"  = "0.42%"
    "D3"  = "38.59"
    "E3"  = "7.79%"
    "D4"  = "5.082"
    "E4"  = "0.70%"
    "D5"  = "0.08075"
    "E5"  = "0.37%"
    "D6"  = "1.945"
    "E6"  = "4.59%"
    "D7"  = "4.190"
    "E7"  = "1.67%"
    "D8"  = "7.937"
    "E8"  = "2.10%"
    "D9"  = "0.9300"
    "E9"  = "0.44%"
    "D10" = "0.1463"
    "E10" = "14.60%"
    "E11" = "2.38%"
    "D12" = "0.09013"
    "E12" = "-0.52%"
    "D13" = "0.03494"
    "E13" = "1.55%"
    "D14" = "0.09792"
    "E14" = "-0.85%"
    "D15" = "0.001389"
    "E15" = "-2.13%"
    "D16" = "0.005975"
    "E16" = "-4.86%"
    "D17" = "3.728"
    "E17" = "-3.48%"
    "E18" = "3.06%"
    "D19" = "0.3463"
    "E19" = "1.56%"
    "D20" = "0.1311"
    "E20" = "0.86%"
    "D21" = "4.787"
    "E21" = "-0.72%"
    "E22" = "-0.21%"
    "D23" = "0.04372"
    "E23" = "0.17%"
    "E24" = "0.21%"
    "D25" = "0.004270"
    "E25" = "-11.80%"
    "D27" = "0.0001300"
    "E27" = "0.00%"
    "D39" = "0.02057"
    "E39" = "4.64%"
    "D40" = "0.05045"
    "E40" = "-1.95%"
    "D41" = "0.007462"
    "E41" = "-0.80%"
    "E42" = "0.14%"
    "D43" = "0.1350"
    "E43" = "-0.34%"
    "D44" = "0.002141"
    "E44" = "1.42%"
    "D45" = "0.008942"
    "E45" = "-9.39%"
    "D46" = "0.00006202"
    "E47" = "0.02%"
    "D48" = "0.002807"
    "E49" = "27.96%"
    "E50" = "0.02%"
    "E51" = "0.02%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
